$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updates to existing rows (2-93) per diff ---
$ws.Range("D2").Value = 44235
$ws.Range("J2").Value = 440
$ws.Range("K2").Value = 17000
$ws.Range("L2").Value = 18000
$ws.Range("M2").Value = 17500
$ws.Range("N2").Value = "`$/caja 15 kilos"
$ws.Range("O2").Value = "Provincia de Limarí"
$ws.Range("P2").Value = 1167
$ws.Range("Q2").Value = 15
$ws.Range("D3").Value = 44335
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 440
$ws.Range("K3").Value = 28000
$ws.Range("L3").Value = 29000
$ws.Range("M3").Value = 28500
$ws.Range("P3").Value = 2375
$ws.Range("D4").Value = 44432
$ws.Range("J4").Value = 540
$ws.Range("K4").Value = 35000
$ws.Range("L4").Value = 36000
$ws.Range("M4").Value = 35500
$ws.Range("P4").Value = 2958
$ws.Range("D5").Value = 44428
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 700
$ws.Range("K5").Value = 36000
$ws.Range("L5").Value = 37000
$ws.Range("M5").Value = 36500
$ws.Range("P5").Value = 3042
$ws.Range("D6").Value = 44428
$ws.Range("I6").Value = "Segunda"
$ws.Range("K6").Value = 29000
$ws.Range("L6").Value = 30000
$ws.Range("M6").Value = 29500
$ws.Range("P6").Value = 2458
$ws.Range("D7").Value = 44362
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 400
$ws.Range("K7").Value = 23500
$ws.Range("L7").Value = 24000
$ws.Range("M7").Value = 23750
$ws.Range("P7").Value = 1979
$ws.Range("D8").Value = 44454
$ws.Range("J8").Value = 640
$ws.Range("K8").Value = 38000
$ws.Range("L8").Value = 39000
$ws.Range("M8").Value = 38500
$ws.Range("P8").Value = 3208
$ws.Range("D9").Value = 44202
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 500
$ws.Range("K9").Value = 24000
$ws.Range("L9").Value = 25000
$ws.Range("M9").Value = 24500
$ws.Range("N9").Value = "`$/caja 15 kilos"
$ws.Range("O9").Value = "Provincia de Limarí"
$ws.Range("P9").Value = 1633
$ws.Range("Q9").Value = 15
$ws.Range("D10").Value = 44424
$ws.Range("J10").Value = 560
$ws.Range("K10").Value = 37000
$ws.Range("L10").Value = 38000
$ws.Range("M10").Value = 37500
$ws.Range("P10").Value = 3125
$ws.Range("D11").Value = 44424
$ws.Range("J11").Value = 400
$ws.Range("K11").Value = 30000
$ws.Range("L11").Value = 31000
$ws.Range("M11").Value = 30500
$ws.Range("P11").Value = 2542
$ws.Range("D12").Value = 44398
$ws.Range("J12").Value = 540
$ws.Range("K12").Value = 24000
$ws.Range("L12").Value = 25000
$ws.Range("M12").Value = 24500
$ws.Range("P12").Value = 2042
$ws.Range("D13").Value = 44398
$ws.Range("I13").Value = "Segunda"
$ws.Range("J13").Value = 300
$ws.Range("K13").Value = 18000
$ws.Range("L13").Value = 19000
$ws.Range("M13").Value = 18500
$ws.Range("P13").Value = 1542
$ws.Range("D14").Value = 44421
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 700
$ws.Range("K14").Value = 33000
$ws.Range("L14").Value = 34000
$ws.Range("M14").Value = 33500
$ws.Range("P14").Value = 2792
$ws.Range("D15").Value = 44421
$ws.Range("I15").Value = "Segunda"
$ws.Range("J15").Value = 400
$ws.Range("K15").Value = 28000
$ws.Range("L15").Value = 29000
$ws.Range("M15").Value = 28500
$ws.Range("P15").Value = 2375
$ws.Range("D16").Value = 44379
$ws.Range("H16").Value = "Americana (o)"
$ws.Range("K16").Value = 37000
$ws.Range("L16").Value = 38000
$ws.Range("M16").Value = 37500
$ws.Range("N16").Value = "`$/caja 25 kilos"
$ws.Range("O16").Value = "Provincia de Limarí"
$ws.Range("P16").Value = 1500
$ws.Range("Q16").Value = 25
$ws.Range("D17").Value = 44379
$ws.Range("J17").Value = 560
$ws.Range("K17").Value = 23000
$ws.Range("L17").Value = 24000
$ws.Range("M17").Value = 23500
$ws.Range("P17").Value = 1958
$ws.Range("D18").Value = 44384
$ws.Range("J18").Value = 540
$ws.Range("K18").Value = 23000
$ws.Range("L18").Value = 24000
$ws.Range("M18").Value = 23500
$ws.Range("P18").Value = 1958
$ws.Range("D19").Value = 44179
$ws.Range("J19").Value = 500
$ws.Range("K19").Value = 18000
$ws.Range("L19").Value = 19000
$ws.Range("M19").Value = 18500
$ws.Range("N19").Value = "`$/caja 12 kilos"
$ws.Range("O19").Value = "Región de Arica y Parinacota"
$ws.Range("P19").Value = 1542
$ws.Range("Q19").Value = 12
$ws.Range("D20").Value = 44179
$ws.Range("I20").Value = "Segunda"
$ws.Range("J20").Value = 400
$ws.Range("K20").Value = 11000
$ws.Range("L20").Value = 12000
$ws.Range("M20").Value = 11500
$ws.Range("N20").Value = "`$/caja 12 kilos"
$ws.Range("O20").Value = "Región de Arica y Parinacota"
$ws.Range("P20").Value = 958
$ws.Range("Q20").Value = 12
$ws.Range("D21").Value = 44412
$ws.Range("J21").Value = 600
$ws.Range("K21").Value = 26000
$ws.Range("L21").Value = 27000
$ws.Range("M21").Value = 26500
$ws.Range("P21").Value = 2208
$ws.Range("D22").Value = 44242
$ws.Range("H22").Value = "Americana (o)"
$ws.Range("J22").Value = 480
$ws.Range("K22").Value = 16000
$ws.Range("L22").Value = 17000
$ws.Range("M22").Value = 16500
$ws.Range("N22").Value = "`$/caja 15 kilos"
$ws.Range("O22").Value = "Provincia de Limarí"
$ws.Range("P22").Value = 1100
$ws.Range("Q22").Value = 15
$ws.Range("D23").Value = 44258
$ws.Range("H23").Value = "Inferno"
$ws.Range("J23").Value = 600
$ws.Range("K23").Value = 15000
$ws.Range("L23").Value = 16000
$ws.Range("M23").Value = 15500
$ws.Range("P23").Value = 1033
$ws.Range("D24").Value = 44377
$ws.Range("J24").Value = 520
$ws.Range("K24").Value = 23000
$ws.Range("L24").Value = 24000
$ws.Range("M24").Value = 23500
$ws.Range("P24").Value = 1958
$ws.Range("D25").Value = 44377
$ws.Range("I25").Value = "Segunda"
$ws.Range("J25").Value = 340
$ws.Range("K25").Value = 19000
$ws.Range("L25").Value = 20000
$ws.Range("M25").Value = 19500
$ws.Range("P25").Value = 1625
$ws.Range("D26").Value = 44435
$ws.Range("J26").Value = 1840
$ws.Range("K26").Value = 34000
$ws.Range("L26").Value = 36000
$ws.Range("M26").Value = 35120
$ws.Range("N26").Value = "`$/caja 12 kilos"
$ws.Range("O26").Value = "Región de Arica y Parinacota"
$ws.Range("P26").Value = 2927
$ws.Range("Q26").Value = 12
$ws.Range("D27").Value = 44169
$ws.Range("J27").Value = 600
$ws.Range("K27").Value = 17000
$ws.Range("L27").Value = 18000
$ws.Range("M27").Value = 17500
$ws.Range("P27").Value = 1458
$ws.Range("D28").Value = 44445
$ws.Range("J28").Value = 600
$ws.Range("K28").Value = 41000
$ws.Range("L28").Value = 42000
$ws.Range("M28").Value = 41500
$ws.Range("P28").Value = 3458
$ws.Range("D29").Value = 44272
$ws.Range("J29").Value = 600
$ws.Range("K29").Value = 16000
$ws.Range("L29").Value = 17000
$ws.Range("M29").Value = 16500
$ws.Range("N29").Value = "`$/caja 15 kilos"
$ws.Range("O29").Value = "Provincia de Limarí"
$ws.Range("P29").Value = 1100
$ws.Range("Q29").Value = 15
$ws.Range("D30").Value = 44449
$ws.Range("I30").Value = "Primera"
$ws.Range("J30").Value = 500
$ws.Range("K30").Value = 42000
$ws.Range("L30").Value = 43000
$ws.Range("M30").Value = 42500
$ws.Range("P30").Value = 3542
$ws.Range("D31").Value = 44333
$ws.Range("J31").Value = 400
$ws.Range("K31").Value = 28000
$ws.Range("L31").Value = 29000
$ws.Range("M31").Value = 28500
$ws.Range("N31").Value = "`$/caja 12 kilos"
$ws.Range("O31").Value = "Región de Arica y Parinacota"
$ws.Range("P31").Value = 2375
$ws.Range("Q31").Value = 12
$ws.Range("D32").Value = 44355
$ws.Range("J32").Value = 460
$ws.Range("K32").Value = 24500
$ws.Range("L32").Value = 25000
$ws.Range("M32").Value = 24750
$ws.Range("N32").Value = "`$/caja 12 kilos"
$ws.Range("O32").Value = "Región de Arica y Parinacota"
$ws.Range("P32").Value = 2062
$ws.Range("Q32").Value = 12
$ws.Range("D33").Value = 44446
$ws.Range("J33").Value = 520
$ws.Range("K33").Value = 41000
$ws.Range("L33").Value = 42000
$ws.Range("M33").Value = 41500
$ws.Range("P33").Value = 3458
$ws.Range("D34").Value = 44391
$ws.Range("J34").Value = 540
$ws.Range("K34").Value = 25000
$ws.Range("L34").Value = 26000
$ws.Range("M34").Value = 25500
$ws.Range("P34").Value = 2125
$ws.Range("D35").Value = 44417
$ws.Range("K35").Value = 28000
$ws.Range("L35").Value = 29000
$ws.Range("M35").Value = 28500
$ws.Range("N35").Value = "`$/caja 12 kilos"
$ws.Range("O35").Value = "Región de Arica y Parinacota"
$ws.Range("P35").Value = 2375
$ws.Range("Q35").Value = 12
$ws.Range("D36").Value = 44417
$ws.Range("I36").Value = "Segunda"
$ws.Range("J36").Value = 400
$ws.Range("K36").Value = 22000
$ws.Range("L36").Value = 23000
$ws.Range("M36").Value = 22500
$ws.Range("P36").Value = 1875
$ws.Range("D37").Value = 44253
$ws.Range("J37").Value = 700
$ws.Range("K37").Value = 15000
$ws.Range("L37").Value = 16000
$ws.Range("M37").Value = 15500
$ws.Range("N37").Value = "`$/caja 15 kilos"
$ws.Range("O37").Value = "Provincia de Limarí"
$ws.Range("P37").Value = 1033
$ws.Range("Q37").Value = 15
$ws.Range("D38").Value = 44342
$ws.Range("J38").Value = 460
$ws.Range("K38").Value = 28000
$ws.Range("L38").Value = 29000
$ws.Range("M38").Value = 28500
$ws.Range("P38").Value = 2375
$ws.Range("D39").Value = 44452
$ws.Range("J39").Value = 600
$ws.Range("K39").Value = 40000
$ws.Range("L39").Value = 41000
$ws.Range("M39").Value = 40500
$ws.Range("P39").Value = 3375
$ws.Range("D40").Value = 44405
$ws.Range("J40").Value = 600
$ws.Range("K40").Value = 26000
$ws.Range("L40").Value = 27000
$ws.Range("M40").Value = 26500
$ws.Range("P40").Value = 2208
$ws.Range("D41").Value = 44246
$ws.Range("J41").Value = 640
$ws.Range("K41").Value = 16000
$ws.Range("L41").Value = 17000
$ws.Range("M41").Value = 16500
$ws.Range("P41").Value = 1100
$ws.Range("D42").Value = 44407
$ws.Range("J42").Value = 760
$ws.Range("K42").Value = 26500
$ws.Range("L42").Value = 27000
$ws.Range("M42").Value = 26750
$ws.Range("P42").Value = 2229
$ws.Range("D43").Value = 44341
$ws.Range("J43").Value = 450
$ws.Range("K43").Value = 28000
$ws.Range("L43").Value = 29000
$ws.Range("M43").Value = 28500
$ws.Range("P43").Value = 2375
$ws.Range("D44").Value = 44274
$ws.Range("J44").Value = 660
$ws.Range("K44").Value = 16000
$ws.Range("L44").Value = 17000
$ws.Range("M44").Value = 16500
$ws.Range("P44").Value = 1100
$ws.Range("D45").Value = 44426
$ws.Range("J45").Value = 600
$ws.Range("K45").Value = 36000
$ws.Range("L45").Value = 37000
$ws.Range("M45").Value = 36500
$ws.Range("P45").Value = 3042
$ws.Range("D46").Value = 44426
$ws.Range("K46").Value = 29000
$ws.Range("L46").Value = 30000
$ws.Range("M46").Value = 29500
$ws.Range("N46").Value = "`$/caja 15 kilos"
$ws.Range("O46").Value = "Provincia de Limarí"
$ws.Range("P46").Value = 1967
$ws.Range("Q46").Value = 15
$ws.Range("D47").Value = 44400
$ws.Range("J47").Value = 700
$ws.Range("K47").Value = 24000
$ws.Range("L47").Value = 25000
$ws.Range("M47").Value = 24500
$ws.Range("P47").Value = 2042
$ws.Range("D48").Value = 44442
$ws.Range("J48").Value = 680
$ws.Range("K48").Value = 38000
$ws.Range("L48").Value = 39000
$ws.Range("M48").Value = 38500
$ws.Range("N48").Value = "`$/caja 12 kilos"
$ws.Range("O48").Value = "Región de Arica y Parinacota"
$ws.Range("P48").Value = 3208
$ws.Range("Q48").Value = 12
$ws.Range("D49").Value = 44453
$ws.Range("K49").Value = 38000
$ws.Range("L49").Value = 39000
$ws.Range("M49").Value = 38500
$ws.Range("N49").Value = "`$/caja 12 kilos"
$ws.Range("O49").Value = "Región de Arica y Parinacota"
$ws.Range("P49").Value = 3208
$ws.Range("Q49").Value = 12
$ws.Range("D50").Value = 44309
$ws.Range("J50").Value = 600
$ws.Range("K50").Value = 20000
$ws.Range("L50").Value = 21000
$ws.Range("M50").Value = 20500
$ws.Range("N50").Value = "`$/caja 15 kilos"
$ws.Range("O50").Value = "Provincia de Limarí"
$ws.Range("P50").Value = 1367
$ws.Range("Q50").Value = 15
$ws.Range("D51").Value = 44370
$ws.Range("K51").Value = 23000
$ws.Range("L51").Value = 24000
$ws.Range("M51").Value = 23500
$ws.Range("N51").Value = "`$/caja 12 kilos"
$ws.Range("O51").Value = "Región de Arica y Parinacota"
$ws.Range("P51").Value = 1958
$ws.Range("Q51").Value = 12
$ws.Range("D52").Value = 44237
$ws.Range("J52").Value = 600
$ws.Range("K52").Value = 17500
$ws.Range("L52").Value = 18000
$ws.Range("M52").Value = 17750
$ws.Range("N52").Value = "`$/caja 15 kilos"
$ws.Range("O52").Value = "Provincia de Limarí"
$ws.Range("P52").Value = 1183
$ws.Range("Q52").Value = 15
$ws.Range("D53").Value = 44167
$ws.Range("J53").Value = 600
$ws.Range("K53").Value = 18000
$ws.Range("L53").Value = 19000
$ws.Range("M53").Value = 18500
$ws.Range("N53").Value = "`$/caja 12 kilos"
$ws.Range("O53").Value = "Región de Arica y Parinacota"
$ws.Range("P53").Value = 1542
$ws.Range("Q53").Value = 12
$ws.Range("D54").Value = 44344
$ws.Range("J54").Value = 540
$ws.Range("K54").Value = 29500
$ws.Range("L54").Value = 30000
$ws.Range("M54").Value = 29750
$ws.Range("P54").Value = 2479
$ws.Range("D55").Value = 44418
$ws.Range("K55").Value = 28000
$ws.Range("L55").Value = 29000
$ws.Range("M55").Value = 28500
$ws.Range("P55").Value = 2375
$ws.Range("D56").Value = 44414
$ws.Range("J56").Value = 700
$ws.Range("K56").Value = 26000
$ws.Range("L56").Value = 27000
$ws.Range("M56").Value = 26500
$ws.Range("P56").Value = 2208
$ws.Range("D57").Value = 44414
$ws.Range("J57").Value = 400
$ws.Range("K57").Value = 25000
$ws.Range("L57").Value = 26000
$ws.Range("M57").Value = 25500
$ws.Range("P57").Value = 1700
$ws.Range("D58").Value = 44165
$ws.Range("J58").Value = 660
$ws.Range("K58").Value = 19000
$ws.Range("L58").Value = 20000
$ws.Range("M58").Value = 19500
$ws.Range("P58").Value = 1625
$ws.Range("D59").Value = 44172
$ws.Range("J59").Value = 560
$ws.Range("K59").Value = 17500
$ws.Range("L59").Value = 18000
$ws.Range("M59").Value = 17750
$ws.Range("P59").Value = 1479
$ws.Range("D60").Value = 44389
$ws.Range("I60").Value = "Primera"
$ws.Range("J60").Value = 600
$ws.Range("K60").Value = 24000
$ws.Range("L60").Value = 25000
$ws.Range("M60").Value = 24500
$ws.Range("N60").Value = "`$/caja 12 kilos"
$ws.Range("O60").Value = "Región de Arica y Parinacota"
$ws.Range("P60").Value = 2042
$ws.Range("Q60").Value = 12
$ws.Range("D61").Value = 44389
$ws.Range("I61").Value = "Segunda"
$ws.Range("J61").Value = 400
$ws.Range("K61").Value = 19000
$ws.Range("L61").Value = 20000
$ws.Range("M61").Value = 19500
$ws.Range("P61").Value = 1625
$ws.Range("D62").Value = 44249
$ws.Range("I62").Value = "Primera"
$ws.Range("J62").Value = 500
$ws.Range("K62").Value = 15000
$ws.Range("L62").Value = 16000
$ws.Range("M62").Value = 15500
$ws.Range("N62").Value = "`$/caja 15 kilos"
$ws.Range("O62").Value = "Provincia de Limarí"
$ws.Range("P62").Value = 1033
$ws.Range("Q62").Value = 15
$ws.Range("D63").Value = 44265
$ws.Range("J63").Value = 600
$ws.Range("K63").Value = 13000
$ws.Range("L63").Value = 14000
$ws.Range("M63").Value = 13500
$ws.Range("N63").Value = "`$/caja 15 kilos"
$ws.Range("O63").Value = "Provincia de Limarí"
$ws.Range("P63").Value = 900
$ws.Range("Q63").Value = 15
$ws.Range("D64").Value = 44447
$ws.Range("I64").Value = "Primera"
$ws.Range("J64").Value = 600
$ws.Range("K64").Value = 42000
$ws.Range("L64").Value = 43000
$ws.Range("M64").Value = 42500
$ws.Range("P64").Value = 3542
$ws.Range("D65").Value = 44260
$ws.Range("J65").Value = 680
$ws.Range("K65").Value = 14000
$ws.Range("L65").Value = 15000
$ws.Range("M65").Value = 14500
$ws.Range("P65").Value = 967
$ws.Range("D66").Value = 44267
$ws.Range("J66").Value = 600
$ws.Range("K66").Value = 13000
$ws.Range("L66").Value = 14000
$ws.Range("M66").Value = 13500
$ws.Range("N66").Value = "`$/caja 15 kilos"
$ws.Range("O66").Value = "Provincia de Limarí"
$ws.Range("P66").Value = 900
$ws.Range("Q66").Value = 15
$ws.Range("D67").Value = 44187
$ws.Range("J67").Value = 520
$ws.Range("K67").Value = 17000
$ws.Range("L67").Value = 18000
$ws.Range("M67").Value = 17500
$ws.Range("P67").Value = 1458
$ws.Range("D68").Value = 44390
$ws.Range("K68").Value = 24000
$ws.Range("L68").Value = 25000
$ws.Range("M68").Value = 24500
$ws.Range("N68").Value = "`$/caja 12 kilos"
$ws.Range("O68").Value = "Región de Arica y Parinacota"
$ws.Range("P68").Value = 2042
$ws.Range("Q68").Value = 12
$ws.Range("D69").Value = 44390
$ws.Range("I69").Value = "Segunda"
$ws.Range("J69").Value = 300
$ws.Range("K69").Value = 19000
$ws.Range("L69").Value = 20000
$ws.Range("M69").Value = 19500
$ws.Range("N69").Value = "`$/caja 12 kilos"
$ws.Range("O69").Value = "Región de Arica y Parinacota"
$ws.Range("P69").Value = 1625
$ws.Range("Q69").Value = 12
$ws.Range("D70").Value = 44386
$ws.Range("K70").Value = 23000
$ws.Range("L70").Value = 24000
$ws.Range("M70").Value = 23500
$ws.Range("P70").Value = 1958
$ws.Range("D71").Value = 44251
$ws.Range("I71").Value = "Primera"
$ws.Range("J71").Value = 600
$ws.Range("K71").Value = 16000
$ws.Range("L71").Value = 17000
$ws.Range("M71").Value = 16500
$ws.Range("N71").Value = "`$/caja 15 kilos"
$ws.Range("O71").Value = "Provincia de Limarí"
$ws.Range("P71").Value = 1100
$ws.Range("Q71").Value = 15
$ws.Range("D72").Value = 44369
$ws.Range("J72").Value = 500
$ws.Range("K72").Value = 24000
$ws.Range("L72").Value = 25000
$ws.Range("M72").Value = 24500
$ws.Range("N72").Value = "`$/caja 12 kilos"
$ws.Range("O72").Value = "Región de Arica y Parinacota"
$ws.Range("P72").Value = 2042
$ws.Range("Q72").Value = 12
$ws.Range("D73").Value = 44433
$ws.Range("J73").Value = 600
$ws.Range("D74").Value = 44397
$ws.Range("J74").Value = 560
$ws.Range("K74").Value = 24500
$ws.Range("M74").Value = 24750
$ws.Range("P74").Value = 2062
$ws.Range("D75").Value = 44363
$ws.Range("J75").Value = 500
$ws.Range("K75").Value = 24000
$ws.Range("L75").Value = 25000
$ws.Range("M75").Value = 24500
$ws.Range("P75").Value = 2042
$ws.Range("D76").Value = 44438
$ws.Range("H76").Value = "Inferno"
$ws.Range("K76").Value = 36000
$ws.Range("L76").Value = 37000
$ws.Range("M76").Value = 36500
$ws.Range("N76").Value = "`$/caja 12 kilos"
$ws.Range("O76").Value = "Región de Arica y Parinacota"
$ws.Range("P76").Value = 3042
$ws.Range("Q76").Value = 12
$ws.Range("D77").Value = 44438
$ws.Range("I77").Value = "Segunda"
$ws.Range("J77").Value = 400
$ws.Range("K77").Value = 30000
$ws.Range("L77").Value = 31000
$ws.Range("M77").Value = 30500
$ws.Range("P77").Value = 2542
$ws.Range("D78").Value = 44372
$ws.Range("K78").Value = 23000
$ws.Range("L78").Value = 24000
$ws.Range("M78").Value = 23500
$ws.Range("P78").Value = 1958
$ws.Range("D79").Value = 44372
$ws.Range("I79").Value = "Segunda"
$ws.Range("J79").Value = 360
$ws.Range("K79").Value = 18000
$ws.Range("L79").Value = 19000
$ws.Range("M79").Value = 18500
$ws.Range("P79").Value = 1542
$ws.Range("D80").Value = 44209
$ws.Range("K80").Value = 20000
$ws.Range("L80").Value = 21000
$ws.Range("M80").Value = 20500
$ws.Range("P80").Value = 1708
$ws.Range("D81").Value = 44356
$ws.Range("K81").Value = 24000
$ws.Range("L81").Value = 25000
$ws.Range("M81").Value = 24500
$ws.Range("P81").Value = 2042
$ws.Range("D82").Value = 44356
$ws.Range("I82").Value = "Segunda"
$ws.Range("J82").Value = 320
$ws.Range("K82").Value = 19000
$ws.Range("L82").Value = 20000
$ws.Range("M82").Value = 19500
$ws.Range("P82").Value = 1625
$ws.Range("D83").Value = 44365
$ws.Range("I83").Value = "Primera"
$ws.Range("J83").Value = 560
$ws.Range("K83").Value = 24000
$ws.Range("L83").Value = 25000
$ws.Range("M83").Value = 24500
$ws.Range("P83").Value = 2042
$ws.Range("D84").Value = 44162
$ws.Range("J84").Value = 700
$ws.Range("K84").Value = 20000
$ws.Range("L84").Value = 21000
$ws.Range("M84").Value = 20500
$ws.Range("N84").Value = "`$/caja 12 kilos"
$ws.Range("O84").Value = "Región de Arica y Parinacota"
$ws.Range("P84").Value = 1708
$ws.Range("Q84").Value = 12
$ws.Range("D85").Value = 44162
$ws.Range("J85").Value = 400
$ws.Range("K85").Value = 34000
$ws.Range("L85").Value = 35000
$ws.Range("M85").Value = 34500
$ws.Range("N85").Value = "`$/caja 25 kilos"
$ws.Range("O85").Value = "Provincia de Limarí"
$ws.Range("P85").Value = 1380
$ws.Range("Q85").Value = 25
$ws.Range("D86").Value = 44410
$ws.Range("I86").Value = "Primera"
$ws.Range("J86").Value = 600
$ws.Range("K86").Value = 26000
$ws.Range("L86").Value = 27000
$ws.Range("M86").Value = 26500
$ws.Range("P86").Value = 2208
$ws.Range("D87").Value = 44411
$ws.Range("J87").Value = 520
$ws.Range("K87").Value = 26000
$ws.Range("L87").Value = 27000
$ws.Range("M87").Value = 26500
$ws.Range("P87").Value = 2208
$ws.Range("D88").Value = 44244
$ws.Range("J88").Value = 600
$ws.Range("K88").Value = 16000
$ws.Range("L88").Value = 17000
$ws.Range("M88").Value = 16500
$ws.Range("N88").Value = "`$/caja 15 kilos"
$ws.Range("O88").Value = "Provincia de Limarí"
$ws.Range("P88").Value = 1100
$ws.Range("Q88").Value = 15
$ws.Range("D89").Value = 44239
$ws.Range("J89").Value = 600
$ws.Range("K89").Value = 17000
$ws.Range("L89").Value = 18000
$ws.Range("M89").Value = 17500
$ws.Range("N89").Value = "`$/caja 15 kilos"
$ws.Range("O89").Value = "Provincia de Limarí"
$ws.Range("P89").Value = 1167
$ws.Range("Q89").Value = 15
$ws.Range("I90").Value = "Primera"
$ws.Range("J90").Value = 520
$ws.Range("K90").Value = 23000
$ws.Range("L90").Value = 24000
$ws.Range("M90").Value = 23500
$ws.Range("P90").Value = 1958
$ws.Range("D91").Value = 44376
$ws.Range("I91").Value = "Segunda"
$ws.Range("J91").Value = 400
$ws.Range("K91").Value = 19000
$ws.Range("L91").Value = 20000
$ws.Range("M91").Value = 19500
$ws.Range("P91").Value = 1625
$ws.Range("D92").Value = 44425
$ws.Range("J92").Value = 540
$ws.Range("K92").Value = 37000
$ws.Range("L92").Value = 38000
$ws.Range("M92").Value = 37500
$ws.Range("P92").Value = 3125
$ws.Range("D93").Value = 44425
$ws.Range("I93").Value = "Segunda"
$ws.Range("J93").Value = 300
$ws.Range("K93").Value = 30000
$ws.Range("L93").Value = 31000
$ws.Range("M93").Value = 30500
$ws.Range("P93").Value = 2542

# --- New row 94 (copy of original row 21 data before its own edit) ---
$ws.Range("A94").Value = 8
$ws.Range("B94").Value = "Terminal La Palmera de La Serena"
$ws.Range("C94").Value = "Coquimbo"
$ws.Range("D94").Value = 44323
$ws.Range("E94").Value = 4
$ws.Range("F94").Value = 100112021
$ws.Range("G94").Value = "Ají"
$ws.Range("H94").Value = "Inferno"
$ws.Range("I94").Value = "Primera"
$ws.Range("J94").Value = 500
$ws.Range("K94").Value = 19000
$ws.Range("L94").Value = 20000
$ws.Range("M94").Value = 19500
$ws.Range("N94").Value = "`$/caja 12 kilos"
$ws.Range("O94").Value = "Región de Arica y Parinacota"
$ws.Range("P94").Value = 1625
$ws.Range("Q94").Value = 12
$ws.Range("R94").Value = "Hortaliza"

# Match the date-style (custom datetime number format) used by column D on other rows
$ws.Range("D94").NumberFormat = $ws.Range("D93").NumberFormat

Write-Host "Edit complete"
